$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bullet char U+2022 used to build the new option-list strings below
$bullet = [char]0x2022

# B2: "Who do you turn to for advice..." answer options
$ws.Range("B2").Value = " $bullet  Pediatrician / doctor   `n $bullet  Childcare provider / teacher   `n $bullet  Home visitor   `n $bullet  Faith leader   `n $bullet  My friends / family   `n $bullet  Others in my town / city / community (e.g., word of mouth, acquaintances, colleagues, neighbors)   `n $bullet  I search online (i.e., Google)   `n $bullet  An online community (e.g., a blog, Facebook group)   `n $bullet  Not Listed (please specify)   "

# B3:B9 share one string: "Never / Sometimes / Frequently" scale used for
# each social-media-platform question row (TikTok, Instagram, YouTube, ...)
$ws.Range("B3:B9").Value = "$bullet Never                                                  $bullet Sometimes                                      $bullet Frequently "

# B10: "What are your favorite places..." answer options
$ws.Range("B10").Value = " $bullet        Public parks   `n $bullet        National / state parks   `n $bullet        Religious buildings (e.g., Mosque, chapel, temple, church)   `n $bullet        Public libraries   `n $bullet        Malls   `n $bullet        Amusement parks   `n $bullet        Bodies of water (e.g., lake, ocean, reservoir)   `n $bullet        My friends' / family's houses   `n $bullet        Childcare programs / school   `n $bullet        My children's non-school programs / classes   `n $bullet        Museums   `n $bullet        Restaurants   `n $bullet        Not Listed (please specify)   "

# B11: "Where are your least favorite places..." answer options
$ws.Range("B11").Value = " $bullet`tPublic parks   `n $bullet`tNational / state parks   `n $bullet`tReligious buildings (e.g., Mosque, chapel, temple, church, etc.)   `n $bullet`tPublic libraries   `n $bullet`tMalls   `n $bullet`tAmusement parks   `n $bullet`tBodies of water (e.g., lake, ocean, reservoir)   `n $bullet`tMy friends' / family's houses   `n $bullet`tChildcare programs / school   `n $bullet`tMy children's non-school programs / classes   `n $bullet`tMuseums   `n $bullet`tRestaurants   `n $bullet`tNot Listed (please specify)   "

# Row-height updates to fit the re-wrapped text
$ws.Rows.Item(2).RowHeight = 196.5
$ws.Rows.Item(10).RowHeight = 270.75

# Row 11's height is unchanged by this revision; pin it back since the
# B11 value edit above causes an autofit recalculation of its wrap height.
$ws.Rows.Item(11).RowHeight = 102.0

